$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'wrestling pants men',
    'mens sliding pants',
    'knee pads for baseball',
    'sliding pad for baseball',
    'compression pants padded knees',
    'goalkeeper knee pads',
    'tights for basketball youth',
    'honeycomb basketball knee pads',
    'knee pad shorts',
    'mens compression capri leggings',
    'boys compression pants with knee pads',
    'rodilleras basketball',
    'honeycomb knee pads basketball',
    'pad pants',
    'long basketball knee pads',
    'men basketball knee',
    'mens baseball sliding short',
    'knee pad for basketball',
    'hex knee pads',
    'construction pants with knee pads',
    'men volleyball knee pads',
    'tight for boys basketball',
    'bjj kneepads',
    'boys padded compression pants',
    'basketball knee tights',
    'youth basketball compression leggings with knee pads',
    'basketball knee compression',
    'youth basketball compression tights',
    'basketball hex knee pads',
    'compression basketball leggings',
    'knee pads for softball',
    'the best knee pads',
    'youth basketball knee pad',
    'mens yoga pants compression',
    'black knee pads basketball',
    'boys basketball knee pads youth',
    'tights basketball',
    'men capri tights',
    'nike basketball knee pads',
    'padded yoga pants',
    'padded knee pants',
    'extra padded knee pads',
    'basketball leg tights',
    'catchers knee pads',
    'basketball tights for youth boys',
    'compression shorts with padding basketball',
    'wrestling pants youth',
    'honeycomb knee pads',
    'baseball leggings for men',
    'pants men basketball',
    'long basketball knee pads adult',
    'yellow leggings for men',
    'mens compression capri',
    'wrestling tights youth boy',
    'basketball youth tights',
    'need pads for basketball',
    'compression capris for men',
    'knee pads for youth basketball',
    'knee pads under pants',
    'goalkeeper padded pants',
    'baseball knee pants',
    'girls softball leggings',
    'youth basketball tights',
    'sliding pad',
    'work pants with knee pads for men',
    'multicam combat pants with knee pads',
    'knee pad honeycomb',
    'cycling leggings men',
    'boys basketball compression tights',
    'compression knee pads work',
    'softball leggings',
    'basketball compression tights',
    'boys basketball pads',
    'basketball knee pads youth boys black',
    'lacrosse knee pads',
    'knee pads youth basketball',
    'hex kneepads',
    'knee basketball',
    'basketball compression pants',
    'knee pad basketball youth',
    'weightlifting pants men',
    'mens knee baseball pants',
    'thick yoga pad',
    'padded basketball compression shorts',
    'under pant knee pads',
    'youth boys basketball knee pads',
    'basketball pad',
    'basketball tights for boys youth',
    'ready man',
    'youth basketball knee pads small',
    'padded compression pants men',
    'long sliding shorts baseball',
    'basketball padded compression shorts men',
    'compression capris men',
    'gym men leggings',
    'sliding shorts youth softball',
    'mens compression pants basketball',
    'softball tights',
    'padded basketball shorts',
    'knee pads basketball men'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
